# Export latest mex files to xlsx. Update segments and articles_db.
# Append 10 new coded-segment rows (411-420) to Sheet1, mirroring the
# formatting of the last existing row (410).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map column letters -> column numbers for A..M
$colNum = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13 }

# Columns whose cell style already carries a Text ("@") number format, so
# assigning a numeric-looking literal stays literal text automatically.
$textFormatCols = @('A','C','E','F','G','L','M')

# Writes $val into $cell as a literal string (never coerced to a number/date)
# while preserving the cell's current style index, by round-tripping through
# a quoted formula and then collapsing it back down to a plain cached value.
function Set-ForcedText($cell, [string]$val) {
    $escaped = $val.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$rows = @(
  @{ row=411; A='●'; B=''; C=''; D='9010'; E='Bacteria:Binomial (genus species)'; F='1: 2617'; G='1: 2637'; H=0; I='Klebsiella pneumoniae'; J=21; K=0.083287; L='Sonia'; M='11/8/18 14:35:00'; ht=16 }
  @{ row=412; A='●'; B=''; C=''; D='21726'; E='Bacteria:Binomial (genus species)'; F='1: 2873'; G='1: 2900'; H=0; I="ycobacterium tubercu- `nlosis"; J=28; K=0.094362; L='Sonia'; M='11/8/18 14:35:00'; ht=30 }
  @{ row=413; A='●'; B=''; C=''; D='21726'; E='Event month'; F='3: 2050'; G='3: 2056'; H=0; I='October'; J=7; K=0.02359; L='Sonia'; M='11/12/18 14:04:00'; ht=16 }
  @{ row=414; A='●'; B=''; C=''; D='21726'; E='Event month'; F='3: 2066'; G='3: 2072'; H=0; I='October'; J=7; K=0.02359; L='Sonia'; M='11/12/18 14:04:00'; ht=16 }
  @{ row=415; A='●'; B=''; C=''; D='21726'; E='Event year'; F='3: 2058'; G='3: 2061'; H=0; I='2006'; J=4; K=0.01348; L='Sonia'; M='11/12/18 14:04:00'; ht=16 }
  @{ row=416; A='●'; B=''; C=''; D='21726'; E='Event year'; F='3: 2074'; G='3: 2077'; H=0; I='2008'; J=4; K=0.01348; L='Sonia'; M='11/12/18 14:04:00'; ht=16 }
  @{ row=417; A='●'; B=''; C=''; D='21726'; E='B'; F='3: 2074'; G='3: 2077'; H=0; I='2008'; J=4; K=0.01348; L='Sonia'; M='11/12/18 14:05:00'; ht=16 }
  @{ row=418; A='●'; B=''; C=''; D='21726'; E='B'; F='3: 2066'; G='3: 2072'; H=0; I='October'; J=7; K=0.02359; L='Sonia'; M='11/12/18 14:05:00'; ht=16 }
  @{ row=419; A='●'; B=''; C=''; D='21726'; E='A'; F='3: 2058'; G='3: 2061'; H=0; I='2006'; J=4; K=0.01348; L='Sonia'; M='11/12/18 14:05:00'; ht=16 }
  @{ row=420; A='●'; B=''; C=''; D='21726'; E='A'; F='3: 2050'; G='3: 2056'; H=0; I='October'; J=7; K=0.02359; L='Sonia'; M='11/12/18 14:05:00'; ht=16 }
)

$srcFormat = $ws.Range("A410:M410")

foreach ($r in $rows) {
    $rowIdx = $r.row

    # Clone the formatting (fill/border/font/number-format) of the last
    # existing data row onto the new row first.
    $dstFormat = $ws.Range("A$($rowIdx):M$($rowIdx)")
    $srcFormat.Copy()
    $dstFormat.PasteSpecial(-4122)  # xlPasteFormats

    foreach ($col in @('A','B','C','D','E','F','G','H','I','J','K','L','M')) {
        $cell = $ws.Cells.Item($rowIdx, $colNum[$col])
        $val = $r[$col]

        if ($col -eq 'H' -or $col -eq 'J' -or $col -eq 'K') {
            # Plain numeric columns.
            $cell.Value = $val
        }
        elseif ($val -eq '') {
            # Blank placeholder cell (Document name / Document group columns).
            $cell.Value = ''
        }
        elseif ($textFormatCols -contains $col) {
            # Column already stored as Text number format -- a direct
            # assignment never gets reinterpreted as a number/date.
            $cell.Value = $val
        }
        elseif ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
            # General-format column but the literal looks numeric (e.g. an
            # ID like "9010" or a bare year like "2006") -- force text so it
            # keeps its original style index instead of becoming a number.
            Set-ForcedText $cell $val
        }
        else {
            $cell.Value = $val
        }
    }

    $ws.Rows.Item($rowIdx).RowHeight = $r.ht
}

Write-Host "Appended rows 411-420 to $($ws.Name)"
